# Edit: "Horarios actualizados Linea 141 - 369"
# Applies updated schedule data across the 3 worksheets (LP1912, LP1912-215, 6203-6173)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---- Sheet1 (LP1912): header metadata ----
$ws1.Cells.Item(2,1).Value = "Última actualización: 06:14:19"
$ws1.Cells.Item(3,1).Value = "Total filas: 48"

# Sheet1 (LP1912) - rows 34-53 data block
$ws1.Cells.Item(34,1).Value = "06:14:19"
$ws1.Cells.Item(34,2).Value = "07:01"
$ws1.Cells.Item(34,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(34,4).Value = 47
$ws1.Cells.Item(34,5).Value = "LP1912"
$ws1.Cells.Item(35,1).Value = "05:49:10"
$ws1.Cells.Item(35,2).Value = "07:05"
$ws1.Cells.Item(35,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(35,4).Value = 76
$ws1.Cells.Item(35,5).Value = "LP1912"
$ws1.Cells.Item(36,1).Value = "05:19:24"
$ws1.Cells.Item(36,2).Value = "07:05"
$ws1.Cells.Item(36,3).Value = "15_ABASTO"
$ws1.Cells.Item(36,4).Value = 106
$ws1.Cells.Item(36,5).Value = "LP1912"
$ws1.Cells.Item(37,1).Value = "05:19:24"
$ws1.Cells.Item(37,2).Value = "07:07"
$ws1.Cells.Item(37,3).Value = "225_GOMEZ"
$ws1.Cells.Item(37,4).Value = 108
$ws1.Cells.Item(37,5).Value = "LP1912"
$ws1.Cells.Item(38,1).Value = "05:19:24"
$ws1.Cells.Item(38,2).Value = "07:11"
$ws1.Cells.Item(38,3).Value = "215A_EL PATO"
$ws1.Cells.Item(38,4).Value = 112
$ws1.Cells.Item(38,5).Value = "LP1912"
$ws1.Cells.Item(39,1).Value = "05:19:24"
$ws1.Cells.Item(39,2).Value = "07:15"
$ws1.Cells.Item(39,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(39,4).Value = 116
$ws1.Cells.Item(39,5).Value = "LP1912"
$ws1.Cells.Item(40,1).Value = "05:49:10"
$ws1.Cells.Item(40,2).Value = "07:16"
$ws1.Cells.Item(40,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(40,4).Value = 87
$ws1.Cells.Item(40,5).Value = "LP1912"
$ws1.Cells.Item(41,1).Value = "05:49:10"
$ws1.Cells.Item(41,2).Value = "07:21"
$ws1.Cells.Item(41,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(41,4).Value = 92
$ws1.Cells.Item(41,5).Value = "LP1912"
$ws1.Cells.Item(42,1).Value = "06:14:19"
$ws1.Cells.Item(42,2).Value = "07:23"
$ws1.Cells.Item(42,3).Value = "10_OLMOS"
$ws1.Cells.Item(42,4).Value = 69
$ws1.Cells.Item(42,5).Value = "LP1912"
$ws1.Cells.Item(43,1).Value = "06:14:19"
$ws1.Cells.Item(43,2).Value = "07:31"
$ws1.Cells.Item(43,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(43,4).Value = 77
$ws1.Cells.Item(43,5).Value = "LP1912"
$ws1.Cells.Item(44,1).Value = "05:49:10"
$ws1.Cells.Item(44,2).Value = "07:32"
$ws1.Cells.Item(44,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(44,4).Value = 103
$ws1.Cells.Item(44,5).Value = "LP1912"
$ws1.Cells.Item(45,1).Value = "05:49:10"
$ws1.Cells.Item(45,2).Value = "07:32"
$ws1.Cells.Item(45,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(45,4).Value = 103
$ws1.Cells.Item(45,5).Value = "LP1912"
$ws1.Cells.Item(46,1).Value = "05:49:10"
$ws1.Cells.Item(46,2).Value = "07:32"
$ws1.Cells.Item(46,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(46,4).Value = 103
$ws1.Cells.Item(46,5).Value = "LP1912"
$ws1.Cells.Item(47,1).Value = "05:49:10"
$ws1.Cells.Item(47,2).Value = "07:37"
$ws1.Cells.Item(47,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(47,4).Value = 108
$ws1.Cells.Item(47,5).Value = "LP1912"
$ws1.Cells.Item(48,1).Value = "05:49:10"
$ws1.Cells.Item(48,2).Value = "07:39"
$ws1.Cells.Item(48,3).Value = "10_OLMOS"
$ws1.Cells.Item(48,4).Value = 110
$ws1.Cells.Item(48,5).Value = "LP1912"
$ws1.Cells.Item(49,1).Value = "06:14:19"
$ws1.Cells.Item(49,2).Value = "07:47"
$ws1.Cells.Item(49,3).Value = "14_ABASTO"
$ws1.Cells.Item(49,4).Value = 93
$ws1.Cells.Item(49,5).Value = "LP1912"
$ws1.Cells.Item(50,1).Value = "05:49:10"
$ws1.Cells.Item(50,2).Value = "07:48"
$ws1.Cells.Item(50,3).Value = "14_ABASTO"
$ws1.Cells.Item(50,4).Value = 119
$ws1.Cells.Item(50,5).Value = "LP1912"
$ws1.Cells.Item(51,1).Value = "06:14:19"
$ws1.Cells.Item(51,2).Value = "07:51"
$ws1.Cells.Item(51,3).Value = "215D_EL PATO"
$ws1.Cells.Item(51,4).Value = 97
$ws1.Cells.Item(51,5).Value = "LP1912"
$ws1.Cells.Item(52,1).Value = "06:14:19"
$ws1.Cells.Item(52,2).Value = "08:00"
$ws1.Cells.Item(52,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(52,4).Value = 106
$ws1.Cells.Item(52,5).Value = "LP1912"
$ws1.Cells.Item(53,1).Value = "06:14:19"
$ws1.Cells.Item(53,2).Value = "08:12"
$ws1.Cells.Item(53,3).Value = "15_ABASTO"
$ws1.Cells.Item(53,4).Value = 118
$ws1.Cells.Item(53,5).Value = "LP1912"

# ---- Sheet2 (LP1912-215): header metadata + new row 14 ----
$ws2.Cells.Item(2,1).Value = "Última actualización: 06:14:19"
$ws2.Cells.Item(3,1).Value = "Total filas: 9"

$ws2.Cells.Item(14,1).Value = "06:14:19"
$ws2.Cells.Item(14,2).Value = "07:51"
$ws2.Cells.Item(14,3).Value = "215D_EL PATO"
$ws2.Cells.Item(14,4).Value = 97
$ws2.Cells.Item(14,5).Value = "LP1912"

# ---- Sheet3 (6203-6173): header metadata + new row 14 ----
$ws3.Cells.Item(2,1).Value = "Última actualización: 06:14:19"
$ws3.Cells.Item(3,1).Value = "Total filas: 9"

$ws3.Cells.Item(14,1).Value = "06:14:19"
$ws3.Cells.Item(14,2).Value = "08:07"
$ws3.Cells.Item(14,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(14,4).Value = 113
$ws3.Cells.Item(14,5).Value = "L6203"
